# Fix the typo "投藍" -> "投籃" ("shoot a basket") in cell A3, and move
# the active selection from G5 to A3 (matching the author's cursor
# position after making the correction).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "投籃"
$ws.Range("A3").Select()
